$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-adjusted timestamp in row 4 (column A)
$ws.Cells.Item(4, 1).Value = 45878.1251727662

# Append the new row 5 with sensor readings
$ws.Cells.Item(5, 1).Value = 45878.16685479179
$ws.Cells.Item(5, 2).Value = 2025
$ws.Cells.Item(5, 3).Value = 37
$ws.Cells.Item(5, 4).Value = 13.17
$ws.Cells.Item(5, 5).Value = 92.28
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 3.14
$ws.Cells.Item(5, 8).Value = "E"
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = "04:00:16"

# Match the style/number format used by the date column for the new row
$ws.Cells.Item(5, 1).NumberFormat = $ws.Cells.Item(4, 1).NumberFormat
